# Add two new columns, I ("I0") and J ("IF"), to the stats sheet.
# I0 is always 1; IF mirrors the existing IP column (H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold font + border + centered alignment)
# from the existing IP header (H1) onto the two new header cells, then
# set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill the data rows: I is a constant 1, J duplicates column H (IP).
$lastRow = 20
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value2
}
